# Reverse the comma-separated "Recorded By" (column G) tokens for affected
# rows, e.g.:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Only these exact "Recorded By" strings are affected by the sync - every
# occurrence of each gets its comma-separated tokens reversed. Other values
# that happen to contain a comma (e.g. "System, admin@admin.com" or
# "dnasr281@gmail.com, admin@admin.com") are intentionally left untouched,
# matching the source data.
$targets = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com" = "backup@backdoor.com, System";
    "System, backup@backdoor.com, system" = "system, backup@backdoor.com, System";
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not $targets.ContainsKey($val)) { continue }

    $cell.Value = $targets[$val]
}
